$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 24
$ws.Range("B24").Value = "Tuesday, Jan 10"
$ws.Range("C24").Value = "8:46 AM"
$ws.Range("D24").Value = "UNKNOWN"
$ws.Range("E24").Value = "Chisinau"
$ws.Range("F24").Value = "(KIV)"
$ws.Range("G24").Value = "SkyUp Airlines "
$ws.Range("H24").Value = "B738"
$ws.Range("I24").Value = "(UR-SQH)"
$ws.Range("J24").Value = "8:44 AM"
$ws.Range("K2").Copy($ws.Range("K24"))
$ws.Range("L24").Value = "0 hours, -2 minutes"
$ws.Range("M2").Copy($ws.Range("M24"))

$ws.Range("A25").Value = 25
$ws.Range("B25").Value = "Tuesday, Jan 10"
$ws.Range("C25").Value = "9:40 AM"
$ws.Range("D25").Value = "FR2136"
$ws.Range("E25").Value = "London"
$ws.Range("F25").Value = "(STN)"
$ws.Range("G25").Value = "Ryanair "
$ws.Range("H25").Value = "B738"
$ws.Range("I25").Value = "(EI-DHE)"
$ws.Range("J25").Value = "9:31 AM"
$ws.Range("K2").Copy($ws.Range("K25"))
$ws.Range("L25").Value = "0 hours, -9 minutes"
$ws.Range("M2").Copy($ws.Range("M25"))

$ws.Range("A26").Value = 26
$ws.Range("B26").Value = "Tuesday, Jan 10"
$ws.Range("C26").Value = "10:15 AM"
$ws.Range("D26").Value = "SK7181"
$ws.Range("E26").Value = "Oslo"
$ws.Range("F26").Value = "(OSL)"
$ws.Range("G26").Value = "SAS "
$ws.Range("H26").Value = "B737"
$ws.Range("I26").Value = "(LN-RPJ)"
$ws.Range("J26").Value = "10:05 AM"
$ws.Range("K2").Copy($ws.Range("K26"))
$ws.Range("L26").Value = "0 hours, -10 minutes"
$ws.Range("M2").Copy($ws.Range("M26"))

$ws.Range("A27").Value = 27
$ws.Range("B27").Value = "Tuesday, Jan 10"
$ws.Range("C27").Value = "10:36 AM"
$ws.Range("D27").Value = "UNKNOWN"
$ws.Range("E27").Value = "Nice"
$ws.Range("F27").Value = "(NCE)"
$ws.Range("G27").Value = "Air X Charter "
$ws.Range("H27").Value = "E35L"
$ws.Range("I27").Value = "(9H-JPC)"
$ws.Range("J27").Value = "10:59 AM"
$ws.Range("K2").Copy($ws.Range("K27"))
$ws.Range("L27").Value = "0 hours, 23 minutes"
$ws.Range("M2").Copy($ws.Range("M27"))

$ws.Range("A28").Value = 28
$ws.Range("B28").Value = "Tuesday, Jan 10"
$ws.Range("C28").Value = "11:15 AM"
$ws.Range("D28").Value = "W61649"
$ws.Range("E28").Value = "Eindhoven"
$ws.Range("F28").Value = "(EIN)"
$ws.Range("G28").Value = "Wizz Air "
$ws.Range("H28").Value = "A320"
$ws.Range("I28").Value = "(HA-LYH)"
$ws.Range("J28").Value = "10:52 AM"
$ws.Range("K2").Copy($ws.Range("K28"))
$ws.Range("L28").Value = "0 hours, -23 minutes"
$ws.Range("M2").Copy($ws.Range("M28"))

$ws.Range("A29").Value = 29
$ws.Range("B29").Value = "Tuesday, Jan 10"
$ws.Range("C29").Value = "11:25 AM"
$ws.Range("D29").Value = "LO3809"
$ws.Range("E29").Value = "Warsaw"
$ws.Range("F29").Value = "(WAW)"
$ws.Range("G29").Value = "LOT "
$ws.Range("H29").Value = "E75S"
$ws.Range("I29").Value = "(SP-LIK)"
$ws.Range("J29").Value = "11:12 AM"
$ws.Range("K2").Copy($ws.Range("K29"))
$ws.Range("L29").Value = "0 hours, -13 minutes"
$ws.Range("M2").Copy($ws.Range("M29"))

$ws.Range("A30").Value = 30
$ws.Range("B30").Value = "Tuesday, Jan 10"
$ws.Range("C30").Value = "2:15 PM"
$ws.Range("D30").Value = "LO3801"
$ws.Range("E30").Value = "Warsaw"
$ws.Range("F30").Value = "(WAW)"
$ws.Range("G30").Value = "LOT "
$ws.Range("H30").Value = "E190"
$ws.Range("I30").Value = "(SP-LMD)"
$ws.Range("J30").Value = "2:09 PM"
$ws.Range("K2").Copy($ws.Range("K30"))
$ws.Range("L30").Value = "0 hours, -6 minutes"
$ws.Range("M2").Copy($ws.Range("M30"))

$ws.Range("A31").Value = 31
$ws.Range("B31").Value = "Tuesday, Jan 10"
$ws.Range("C31").Value = "4:05 PM"
$ws.Range("D31").Value = "LO3807"
$ws.Range("E31").Value = "Warsaw"
$ws.Range("F31").Value = "(WAW)"
$ws.Range("G31").Value = "LOT "
$ws.Range("H31").Value = "E170"
$ws.Range("I31").Value = "(SP-LDH)"
$ws.Range("J31").Value = "3:56 PM"
$ws.Range("K2").Copy($ws.Range("K31"))
$ws.Range("L31").Value = "0 hours, -9 minutes"
$ws.Range("M2").Copy($ws.Range("M31"))

$ws.Range("A32").Value = 32
$ws.Range("B32").Value = "Tuesday, Jan 10"
$ws.Range("C32").Value = "5:55 PM"
$ws.Range("D32").Value = "FR3472"
$ws.Range("E32").Value = "London"
$ws.Range("F32").Value = "(LTN)"
$ws.Range("G32").Value = "Ryanair "
$ws.Range("H32").Value = "B738"
$ws.Range("I32").Value = "(EI-EMR)"
$ws.Range("J32").Value = "5:58 PM"
$ws.Range("K2").Copy($ws.Range("K32"))
$ws.Range("L32").Value = "0 hours, 3 minutes"
$ws.Range("M2").Copy($ws.Range("M32"))

$ws.Range("A33").Value = 33
$ws.Range("B33").Value = "Tuesday, Jan 10"
$ws.Range("C33").Value = "7:50 PM"
$ws.Range("D33").Value = "RK3202"
$ws.Range("E33").Value = "Manchester"
$ws.Range("F33").Value = "(MAN)"
$ws.Range("G33").Value = "Ryanair "
$ws.Range("H33").Value = "B738"
$ws.Range("I33").Value = "(G-RUKH)"
$ws.Range("J33").Value = "7:36 PM"
$ws.Range("K2").Copy($ws.Range("K33"))
$ws.Range("L33").Value = "0 hours, -14 minutes"
$ws.Range("M2").Copy($ws.Range("M33"))

$ws.Range("A34").Value = 34
$ws.Range("B34").Value = "Tuesday, Jan 10"
$ws.Range("C34").Value = "9:15 PM"
$ws.Range("D34").Value = "W95153"
$ws.Range("E34").Value = "London"
$ws.Range("F34").Value = "(LTN)"
$ws.Range("G34").Value = "Wizz Air "
$ws.Range("H34").Value = "A321"
$ws.Range("I34").Value = "(G-WUKG)"
$ws.Range("J34").Value = "8:54 PM"
$ws.Range("K2").Copy($ws.Range("K34"))
$ws.Range("L34").Value = "0 hours, -21 minutes"
$ws.Range("M2").Copy($ws.Range("M34"))

$ws.Range("A35").Value = 35
$ws.Range("B35").Value = "Tuesday, Jan 10"
$ws.Range("C35").Value = "11:50 PM"
$ws.Range("D35").Value = "LO3803"
$ws.Range("E35").Value = "Warsaw"
$ws.Range("F35").Value = "(WAW)"
$ws.Range("G35").Value = "LOT "
$ws.Range("H35").Value = "E190"
$ws.Range("I35").Value = "(SP-LME)"
$ws.Range("J35").Value = "11:35 PM"
$ws.Range("K2").Copy($ws.Range("K35"))
$ws.Range("L35").Value = "0 hours, -15 minutes"
$ws.Range("M2").Copy($ws.Range("M35"))

$ws.Range("A36").Value = 36
$ws.Range("B36").Value = "Wednesday, Jan 11"
$ws.Range("C36").Value = "12:27 AM"
$ws.Range("D36").Value = "AEG481"
$ws.Range("E36").Value = "Ljubljana"
$ws.Range("F36").Value = "(LJU)"
$ws.Range("G36").Value = "Airest "
$ws.Range("H36").Value = "SF34"
$ws.Range("I36").Value = "(ES-LSI)"
$ws.Range("J36").Value = "12:21 AM"
$ws.Range("K2").Copy($ws.Range("K36"))
$ws.Range("L36").Value = "0 hours, -6 minutes"
$ws.Range("M2").Copy($ws.Range("M36"))

